$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "I handleliste"

# Fill "j" in column F for all data rows (rows 2-29, 31-32); row 30 does not exist in the sheet.
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,31,32)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = "j"
}

# Update the active selection to match the edited workbook state
[void]$ws.Range("G26").Select()
